$wb = $excel.ActiveWorkbook

# Sheet2: rename xpath-based locator headers to "_location" headers,
# and update the Sign-in button locator value to the new XPath expression.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("C1").Value2 = "Singin_nav_button_location"
$ws2.Range("D1").Value2 = "Email_location"
$ws2.Range("F1").Value2 = "Next_button_location"
$ws2.Range("G1").Value2 = "Password_location"
$ws2.Range("I1").Value2 = "SingInBtn_location"
$ws2.Range("C2").Value2 = "(//a[@data-task='signin'])[1]"

# Make Sheet2 the active sheet/tab and select C2 (mirrors the workbook's
# saved view state: activeTab on Sheet2, tabSelected moves off Sheet1).
$ws2.Activate() | Out-Null
$ws2.Range("C2").Select() | Out-Null
